# Auto-applies cell value updates for cryptos.xlsx per commit diff (Tue Jun 20 09:53:53 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.738.82"
$ws.Range("E2").Value = "  +1.36%  "

$ws.Range("D3").Value = "1.726.28"
$ws.Range("E3").Value = "  +0.22%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9975"
$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.20"
$ws.Range("E5").Value = "  -0.91%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9982"
$ws.Range("E6").Value = "  -0.17%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4836"
$ws.Range("E7").Value = "  -0.77%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2579"
$ws.Range("E8").Value = "  -0.37%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06182"
$ws.Range("E9").Value = "  -0.18%  "

$ws.Range("D10").Value = "1.719.38"
$ws.Range("E10").Value = "  -0.19%  "

$ws.Range("E11").Value = "  +2.78%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.06853"
$ws.Range("E12").Value = "  -1.80%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6047"
$ws.Range("E13").Value = "  +1.42%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.457"
$ws.Range("E14").Value = "  -1.62%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.10"
$ws.Range("E15").Value = "  -0.06%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9982"
$ws.Range("E16").Value = "  -0.18%  "

$ws.Range("D17").Value = "26.558.27"
$ws.Range("E17").Value = "  +0.67%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9975"
$ws.Range("E18").Value = "  -0.23%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007152"
$ws.Range("E19").Value = "  -0.67%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.37"
$ws.Range("E20").Value = "  +0.49%  "

$ws.Range("D21").Value = "1.941.43"
$ws.Range("E21").Value = "  -0.18%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.405"
$ws.Range("E22").Value = "  -0.80%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.561"
$ws.Range("E23").Value = "  +0.96%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.054"
$ws.Range("E24").Value = "  -0.86%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "139.69"
$ws.Range("E25").Value = "  +1.35%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.29"
$ws.Range("E26").Value = "  +0.35%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.775"
$ws.Range("E27").Value = "  +3.20%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "106.73"
$ws.Range("E28").Value = "  +0.06%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.376"
$ws.Range("E29").Value = "  -1.64%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.012"
$ws.Range("E30").Value = "  +2.41%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07919"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.661"
$ws.Range("E32").Value = "  -0.03%  "

$ws.Range("E33").Value = "  +0.13%  "

$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.593"
$ws.Range("E34").Value = "  -0.53%  "

$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9994"
$ws.Range("E35").Value = "  +0.42%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6174"
$ws.Range("E36").Value = "  -0.98%  "

$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9348"
$ws.Range("E37").Value = "  +0.45%  "

$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.009"
$ws.Range("E38").Value = "  +2.62%  "

$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.449"
$ws.Range("E39").Value = "  +2.47%  "

$ws.Range("B40").Value = "PaxDollar"
$ws.Range("C40").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9974"
$ws.Range("E40").Value = "  -0.17%  "

$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01496"
$ws.Range("E41").Value = "  +1.66%  "

$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.617"
$ws.Range("E42").Value = "  +3.71%  "

$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.76"
$ws.Range("E43").Value = "  -0.31%  "

$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3825"
$ws.Range("E44").Value = "  -0.23%  "

$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.795"
$ws.Range("E45").Value = "  -0.48%  "

$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1153"
$ws.Range("E46").Value = "  -0.92%  "

$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05361"
$ws.Range("E47").Value = "  -0.07%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.909"
$ws.Range("E48").Value = "  +3.17%  "

$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.04"
$ws.Range("E49").Value = "  -0.07%  "

$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.243"
$ws.Range("E50").Value = "  +1.39%  "

$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.29"
$ws.Range("E51").Value = "  +0.85%  "

